# Update 2018-05-24#7  -  setup api routes part 5
# Adds a new "Devices" error-code section (rows 245-257) to Sheet1, mirroring
# the existing table layout (bold/left-aligned section header row, merged
# A:C on the header, plain data rows below).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---- Section header row (245): bold, left aligned, merged A245:C245 ----
$ws.Range("A245").Value = "Devices"
$headerRange = $ws.Range("A245:C245")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4131  # xlLeft
$headerRange.Merge()

# ---- Data rows (246-257) ----
# NOTE: the write order below intentionally matches the original author's
# entry order (column-by-column for the first sub-table, row-by-row for the
# second) so that new shared-string entries are interned in the same
# sequence as the canonical document (sharedStrings.xml append order).

# -- "SaveDevice" sub-table (rows 246-251): column A, then column B, then column C --
$ws.Range("A246").Value = 2401
$ws.Range("A247").Value = 2402
$ws.Range("A248").Value = 2403
$ws.Range("A249").Value = 2404
$ws.Range("A250").Value = 2405
$ws.Range("A251").Value = 2406

$ws.Range("B246").Value = "Customer Id cannot be null or empty string."
$ws.Range("B247").Value = "Device Type Id not found."
$ws.Range("B248").Value = "Device Name (default) cannot be null or empty string."
$ws.Range("B249").Value = "Customer Id is not found."
$ws.Range("B250").Value = "Device Id is not found."
$ws.Range("B251").Value = "Device Name (default) already exists."

$ws.Range("C246").Value = "SaveDevice"
$ws.Range("C247").Value = "SaveDevice"
$ws.Range("C248").Value = "SaveDevice"
$ws.Range("C249").Value = "SaveDevice"
$ws.Range("C250").Value = "SaveDevice"
$ws.Range("C251").Value = "SaveDevice"

# -- "SaveDeviceML" sub-table (rows 252-257): row by row (A, B, C) --
$ws.Range("A252").Value = 2407
$ws.Range("B252").Value = "Customer Id cannot be null or empty string."
$ws.Range("C252").Value = "SaveDeviceML"

$ws.Range("A253").Value = 2408
$ws.Range("B253").Value = "Lang Id cannot be null or empty string."
$ws.Range("C253").Value = "SaveDeviceML"

$ws.Range("A254").Value = 2409
$ws.Range("B254").Value = "Lang Id not exist."
$ws.Range("C254").Value = "SaveDeviceML"

$ws.Range("A255").Value = 2410
$ws.Range("B255").Value = "Device Id cannot be null or empty string."
$ws.Range("C255").Value = "SaveDeviceML"

$ws.Range("A256").Value = 2411
$ws.Range("B256").Value = "Device Id is not found."
$ws.Range("C256").Value = "SaveDeviceML"

$ws.Range("A257").Value = 2412
$ws.Range("B257").Value = "Device Name (ML) is already exists."
$ws.Range("C257").Value = "SaveDeviceML"

# ---- View state: scroll position + new selection below the added rows ----
$ws.Activate() | Out-Null
$excel.ActiveWindow.ScrollRow = 234
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("A258:XFD260").Select() | Out-Null
